$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("December")

$data = @(
    @(2, 1344, 1185, 159, 'We borrowerd more than we lent', '', '1.13 : 1'),
    @(3, 438, 477, -39, '', 'We lent more than we borrowed', '0.92 : 1'),
    @(4, 1029, 1108, -79, '', 'We lent more than we borrowed', '0.93 : 1'),
    @(5, 61, 120, -59, '', 'We lent more than we borrowed', '0.51 : 1'),
    @(6, 927, 1208, -281, '', 'We lent more than we borrowed', '0.77 : 1'),
    @(7, 105, 164, -59, '', 'We lent more than we borrowed', '0.64 : 1'),
    @(8, 84, 129, -45, '', 'We lent more than we borrowed', '0.65 : 1'),
    @(9, 77, 72, 5, 'We borrowerd more than we lent', '', '1.07 : 1'),
    @(10, 0, 50, -50, '', 'We lent more than we borrowed', '0.00 : 1'),
    @(11, 0, 0, 0, '', '', ''),
    @(12, 15, 30, -15, '', 'We lent more than we borrowed', '0.50 : 1'),
    @(13, 185, 64, 121, 'We borrowerd more than we lent', '', '2.89 : 1'),
    @(14, 119, 199, -80, '', 'We lent more than we borrowed', '0.60 : 1'),
    @(15, 80, 89, -9, '', 'We lent more than we borrowed', '0.90 : 1'),
    @(16, 53, 143, -90, '', 'We lent more than we borrowed', '0.37 : 1'),
    @(17, 574, 484, 90, 'We borrowerd more than we lent', '', '1.19 : 1'),
    @(18, 96, 91, 5, 'We borrowerd more than we lent', '', '1.05 : 1'),
    @(19, 620, 302, 318, 'We borrowerd more than we lent', '', '2.05 : 1'),
    @(20, 0, 53, -53, '', 'We lent more than we borrowed', '0.00 : 1'),
    @(21, 399, 340, 59, 'We borrowerd more than we lent', '', '1.17 : 1'),
    @(22, 46, 119, -73, '', 'We lent more than we borrowed', '0.39 : 1'),
    @(23, 675, 346, 329, 'We borrowerd more than we lent', '', '1.95 : 1'),
    @(24, 1586, 1122, 464, 'We borrowerd more than we lent', '', '1.41 : 1'),
    @(25, 149, 307, -158, '', 'We lent more than we borrowed', '0.49 : 1'),
    @(26, 0, 0, 0, '', '', ''),
    @(27, 216, 187, 29, 'We borrowerd more than we lent', '', '1.16 : 1'),
    @(28, 51, 103, -52, '', 'We lent more than we borrowed', '0.50 : 1'),
    @(29, 430, 418, 12, 'We borrowerd more than we lent', '', '1.03 : 1'),
    @(30, 45, 15, 30, 'We borrowerd more than we lent', '', '3.00 : 1'),
    @(31, 63, 245, -182, '', 'We lent more than we borrowed', '0.26 : 1'),
    @(32, 354, 540, -186, '', 'We lent more than we borrowed', '0.66 : 1'),
    @(33, 326, 496, -170, '', 'We lent more than we borrowed', '0.66 : 1'),
    @(34, 160, 116, 44, 'We borrowerd more than we lent', '', '1.38 : 1'),
    @(35, 780, 871, -91, '', 'We lent more than we borrowed', '0.90 : 1'),
    @(36, 198, 424, -226, '', 'We lent more than we borrowed', '0.47 : 1'),
    @(37, 413, 270, 143, 'We borrowerd more than we lent', '', '1.53 : 1'),
    @(38, 28, 128, -100, '', 'We lent more than we borrowed', '0.22 : 1'),
    @(39, 15, 93, -78, '', 'We lent more than we borrowed', '0.16 : 1'),
    @(40, 41, 139, -98, '', 'We lent more than we borrowed', '0.29 : 1'),
    @(41, 2, 28, -26, '', 'We lent more than we borrowed', '0.07 : 1'),
    @(42, 5, 24, -19, '', 'We lent more than we borrowed', '0.21 : 1'),
    @(43, 0, 0, 0, '', '', ''),
    @(44, 52, 71, -19, '', 'We lent more than we borrowed', '0.73 : 1'),
    @(45, 84, 147, -63, '', 'We lent more than we borrowed', '0.57 : 1'),
    @(46, 511, 563, -52, '', 'We lent more than we borrowed', '0.91 : 1'),
    @(47, 1019, 483, 536, 'We borrowerd more than we lent', '', '2.11 : 1'),
    @(48, 194, 554, -360, '', 'We lent more than we borrowed', '0.35 : 1'),
    @(49, 539, 227, 312, 'We borrowerd more than we lent', '', '2.37 : 1'),
    @(50, 868, 497, 371, 'We borrowerd more than we lent', '', '1.75 : 1'),
    @(51, 193, 169, 24, 'We borrowerd more than we lent', '', '1.14 : 1'),
    @(52, 406, 428, -22, '', 'We lent more than we borrowed', '0.95 : 1'),
    @(53, 107, 226, -119, '', 'We lent more than we borrowed', '0.47 : 1'),
    @(54, 20, 200, -180, '', 'We lent more than we borrowed', '0.10 : 1'),
    @(55, 279, 197, 82, 'We borrowerd more than we lent', '', '1.42 : 1'),
)

foreach ($row in $data) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $d = $row[3]
    $eText = $row[4]
    $fText = $row[5]
    $gText = $row[6]

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d

    if ($eText -ne "") {
        $ws.Cells.Item($r, 5).Value = $eText
    }
    if ($fText -ne "") {
        $ws.Cells.Item($r, 6).Value = $fText
    }
    if ($gText -ne "") {
        $ws.Cells.Item($r, 7).Value = $gText
    }
}
